$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Insert a new column before column D (the "notes" column), shifting it to E
$ws.Columns("D:D").Insert()

# Set the header for new column D
$ws.Range("D1").Value = "surveyor"

# Set surveyor values for each data row
$ws.Range("D2").Value = "Grace"
$ws.Range("D3").Value = "Grace"
$ws.Range("D4").Value = "Erik"
$ws.Range("D5").Value = "Erik"
$ws.Range("D6").Value = "Erik"
$ws.Range("D7").Value = "Grace"
$ws.Range("D8").Value = "Grace"
$ws.Range("D9").Value = "Erik"
$ws.Range("D10").Value = "Erik"
$ws.Range("D11").Value = "Erik"

# Update selection on the sheet
$ws.Range("G30").Select()
